$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3)
$ws.Range("B3").Value = "0.1.7"

# Update Status value (row 6)
$ws.Range("B6").Value = "draft"

# Update Date value (row 8)
$ws.Range("B8").Value = "2024-08-27T12:23:18-05:00"

# Row 10 is "Contact" / "No display for ContactDetail" -> update value to publisher contact
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# Row 11 is already a duplicate "Contact" row with the same old value -> update value to the person contact
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# Insert a new row after row 11 for the Jurisdiction entry, pushing Description/Purpose/Copyright/Immutable down by one
$ws.Rows.Item(12).Insert()

$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""
